# Update the standard-error rows (theta_se row 4, lambda_se row 6) of the
# police-report-only appendix table with the finalized values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: theta_se
$ws.Range("B4").Value = "(3.74)"
$ws.Range("C4").Value = "(1.03)"
$ws.Range("D4").Value = "(2.02)"
$ws.Range("E4").Value = "(1.83)"
$ws.Range("F4").Value = "(1.45)"
$ws.Range("G4").Value = "(1.78)"
$ws.Range("H4").Value = "(2.41)"
$ws.Range("I4").Value = "(2.53)"
$ws.Range("J4").Value = "(2.27)"
$ws.Range("K4").Value = "(2.3)"
$ws.Range("L4").Value = "(2.89)"

# Row 6: lambda_se
$ws.Range("B6").Value = "(3.18)"
$ws.Range("C6").Value = "(0.79)"
$ws.Range("D6").Value = "(1.42)"
$ws.Range("E6").Value = "(1.52)"
$ws.Range("F6").Value = "(0.75)"
$ws.Range("G6").Value = "(1.49)"
$ws.Range("H6").Value = "(2.04)"
$ws.Range("I6").Value = "(2.12)"
$ws.Range("J6").Value = "(1.38)"
$ws.Range("K6").Value = "(1.88)"
$ws.Range("L6").Value = "(2.85)"
